$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# B2 ("Past") and D2 ("Take notes on enemies") are unchanged.
# C2 ("Review Transform Manipulation") keeps its yellow fill (no alignment).
$ws.Range("C2").Interior.ColorIndex = 6

# New cell E2: "Download a Zombie Sprite", centered (no fill, no wrap).
$ws.Range("E2").Value = "Download a Zombie Sprite"
$ws.Range("E2").HorizontalAlignment = -4108

# --- Row 3 (ht 30) ---
# C3 ("Review Component Access") unchanged.
# D3 text updated, and now also gets the yellow fill in addition to its
# existing center+wrap alignment.
$ws.Range("D3").Value = "Pseudocode out Zombie movement behavior"
$ws.Range("D3").Interior.ColorIndex = 6

# New cell E3: "Figure out how to change direction".
$ws.Range("E3").Value = "Figure out how to change direction"

# --- Row 4 ---
# C4 ("Make Plan for the rest of the week") unchanged.
# New cell D4: "Get an object moving", with the plain yellow fill style.
$ws.Range("D4").Value = "Get an object moving"
$ws.Range("D4").Interior.ColorIndex = 6

# New (empty) cell E4, matching the formatting used for E3.
$ws.Range("E4").HorizontalAlignment = -4108

# --- Row 5 (now ht 30) ---
# C5 ("Study GnG enemy (zombie) behavior") keeps its yellow-fill style.
$ws.Range("C5").Interior.ColorIndex = 6

# New cell D5: "Attempt a lerp zombie movement implementation", with the
# yellow fill + center + wrap style (same combo as D3).
$ws.Range("D5").Value = "Attempt a lerp zombie movement implementation"
$ws.Range("D5").Interior.ColorIndex = 6
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").WrapText = $true
$ws.Rows("5:5").RowHeight = 30

# Widen column E now that it holds real content.
$ws.Columns("E:E").ColumnWidth = 24.7109375

# Final selection left on D5, matching the last-edited cell.
$ws.Range("D5").Select()
